$d = $word.ActiveDocument

# 1. Title heading change (first occurrence, used at top of doc)
$d.Content.Find.Execute(
    "Play Im King Free Slot Game Online | Win Significant Payouts", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play IM KING - Free Online Slot Game", 2)

# Pros list reshuffled - process in reverse dependency order to avoid
# collisions between newly-written text and text we still need to find.
# Original order:
#   32: Range of bonus features including Wild symbols and free spins
#   33: Flexible betting range from 0.35 cents up to 350€
#   34: Visually appealing graphics depicting Chinese empire symbols
#   35: Opportunity to win significant payouts
# New order:
#   32: Opportunity to win significant payouts
#   33: Range of bonus features (Wild symbols, free spins)
#   34: Flexible betting range
#   35: Visually appealing graphics

$d.Content.Find.Execute(
    "Opportunity to win significant payouts", $true, $false, $false, $false, $false,
    $true, 1, $false, "Visually appealing graphics", 2)

$d.Content.Find.Execute(
    "Visually appealing graphics depicting Chinese empire symbols", $true, $false, $false, $false, $false,
    $true, 1, $false, "Flexible betting range", 2)

$d.Content.Find.Execute(
    "Flexible betting range from 0.35 cents up to 350€", $true, $false, $false, $false, $false,
    $true, 1, $false, "Range of bonus features (Wild symbols, free spins)", 2)

$d.Content.Find.Execute(
    "Range of bonus features including Wild symbols and free spins", $true, $false, $false, $false, $false,
    $true, 1, $false, "Opportunity to win significant payouts", 2)

# Cons list wording tweaks
$d.Content.Find.Execute(
    "Music can be a bit annoying", $true, $false, $false, $false, $false,
    $true, 1, $false, "Annoying music", 2)

$d.Content.Find.Execute(
    "No information provided on volatility and RTP", $true, $false, $false, $false, $false,
    $true, 1, $false, "No information on volatility and RTP", 2)

# Bold "title" run near the end of the document
$d.Content.Find.Execute(
    "Play Im King Free Slot Game Online | Win Significant Payouts", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play IM KING - Free Online Slot Game", 2)

# Italic meta-description run
$d.Content.Find.Execute(
    "Im King is a 5x3 online slot game with Wild symbols, free spins, and bonuses. Play now for free and win significant payouts. No download or registration required.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Read our review of IM KING, a free online slot game with exciting bonus features.", 2)
